$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the text of the shared "Include " string (currently at row 15, column F)
# before the row insertion shifts it down to row 16.
$ws.Range("F15").Value = "Included "

# Insert a new row at position 6 (shifts existing rows 6.. down by one)
$ws.Rows.Item(6).Insert()

# Populate the new row 6 with values, mirroring row 5's layout
$ws.Range("A6").Value = "輕輪動平衡"
$ws.Range("C6").Value = "再生車修"
$ws.Range("E6").Value = "A"
$ws.Range("H6").Value = "須符合規格"

# Match the saved view/selection state (user ended up with A6 selected)
[void]$ws.Range("A6").Select()

# The hidden _FilterDatabase name still spans the original (pre-insert) range;
# extend it by one row to keep covering the whole table (A1:H69).
$filterName = $wb.Names.Item(1)
$filterName.RefersTo = "=" + $ws.Name + "!`$A`$1:`$H`$69"
